$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2222222222222222
$ws.Range("C2").Value = 0.5095785440613027
$ws.Range("J2").Value = 0.01149425287356322
$ws.Range("P2").Value = 0.1379310344827586
$ws.Range("S2").Value = 0.1187739463601533
# Row 3
$ws.Range("B3").Value = 0.01492537313432836
$ws.Range("C3").Value = 0.02238805970149254
$ws.Range("J3").Value = 0.02238805970149254
$ws.Range("P3").Value = 0.8208955223880597
$ws.Range("S3").Value = 0.1194029850746269
# Row 4
$ws.Range("J4").Value = 0.08888888888888889
$ws.Range("P4").Value = 0.6888888888888889
$ws.Range("S4").Value = 0.2222222222222222
# Row 6
$ws.Range("B6").Value = 0.06532663316582915
$ws.Range("D6").Value = 0.03015075376884422
$ws.Range("F6").Value = 0.06030150753768844
$ws.Range("J6").Value = 0.2361809045226131
$ws.Range("O6").Value = 0.02010050251256281
$ws.Range("Q6").Value = 0.1708542713567839
$ws.Range("R6").Value = 0.07537688442211055
$ws.Range("S6").Value = 0.3417085427135678
# Row 7
$ws.Range("B7").Value = 0.09326424870466321
$ws.Range("D7").Value = 0.02590673575129534
$ws.Range("F7").Value = 0.06217616580310881
$ws.Range("J7").Value = 0.1658031088082902
$ws.Range("O7").Value = 0.02072538860103627
$ws.Range("Q7").Value = 0.1865284974093264
$ws.Range("R7").Value = 0.05699481865284974
$ws.Range("S7").Value = 0.38860103626943
# Row 8
$ws.Range("B8").Value = 0.09767441860465116
$ws.Range("D8").Value = 0.02558139534883721
$ws.Range("E8").Value = 0.002325581395348837
$ws.Range("F8").Value = 0.04651162790697674
$ws.Range("J8").Value = 0.09767441860465116
$ws.Range("O8").Value = 0.0186046511627907
$ws.Range("Q8").Value = 0.1651162790697674
$ws.Range("R8").Value = 0.1023255813953488
$ws.Range("S8").Value = 0.4441860465116279
# Row 9
$ws.Range("B9").Value = 0.0736196319018405
$ws.Range("D9").Value = 0.0245398773006135
$ws.Range("F9").Value = 0.0245398773006135
$ws.Range("J9").Value = 0.1288343558282209
$ws.Range("O9").Value = 0.01840490797546012
$ws.Range("Q9").Value = 0.1840490797546012
$ws.Range("R9").Value = 0.08588957055214724
$ws.Range("S9").Value = 0.4601226993865031
# Row 10
$ws.Range("B10").Value = 0.100088573959256
$ws.Range("D10").Value = 0.01860053144375554
$ws.Range("F10").Value = 0.07883082373782108
$ws.Range("J10").Value = 0.129317980513729
$ws.Range("O10").Value = 0.01328609388839681
$ws.Range("Q10").Value = 0.2214348981399469
$ws.Range("R10").Value = 0.08060230292294066
$ws.Range("S10").Value = 0.3578387953941541
# Row 11
$ws.Range("G11").Value = 0.1423611111111111
$ws.Range("J11").Value = 0.07291666666666667
$ws.Range("K11").Value = 0.1701388888888889
$ws.Range("L11").Value = 0.6006944444444444
$ws.Range("S11").Value = 0.01388888888888889
# Row 12
$ws.Range("G12").Value = 0.7388888888888889
$ws.Range("J12").Value = 0.1888888888888889
$ws.Range("K12").Value = 0.01111111111111111
$ws.Range("L12").Value = 0.03333333333333333
$ws.Range("S12").Value = 0.02777777777777778
# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.1025641025641026
# Row 15
$ws.Range("F15").Value = 0.03414634146341464
$ws.Range("H15").Value = 0.1658536585365854
$ws.Range("I15").Value = 0.07317073170731707
$ws.Range("J15").Value = 0.3560975609756097
$ws.Range("K15").Value = 0.08292682926829269
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.05853658536585366
$ws.Range("S15").Value = 0.2048780487804878
# Row 16
$ws.Range("F16").Value = 0.01169590643274854
$ws.Range("H16").Value = 0.2339181286549707
$ws.Range("I16").Value = 0.1052631578947368
$ws.Range("J16").Value = 0.3567251461988304
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("M16").Value = 0.01754385964912281
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1169590643274854
# Row 17
$ws.Range("F17").Value = 0.02612826603325416
$ws.Range("H17").Value = 0.1923990498812352
$ws.Range("I17").Value = 0.06175771971496437
$ws.Range("J17").Value = 0.4038004750593824
$ws.Range("K17").Value = 0.1211401425178147
$ws.Range("M17").Value = 0.01187648456057007
$ws.Range("O17").Value = 0.0665083135391924
$ws.Range("S17").Value = 0.1163895486935867
# Row 18
$ws.Range("F18").Value = 0.01744186046511628
$ws.Range("H18").Value = 0.1918604651162791
$ws.Range("I18").Value = 0.06976744186046512
$ws.Range("J18").Value = 0.4186046511627907
$ws.Range("K18").Value = 0.09883720930232558
$ws.Range("M18").Value = 0.01162790697674419
$ws.Range("O18").Value = 0.1046511627906977
$ws.Range("S18").Value = 0.0872093023255814
# Row 19
$ws.Range("F19").Value = 0.01704035874439462
$ws.Range("H19").Value = 0.2179372197309417
$ws.Range("I19").Value = 0.08251121076233184
$ws.Range("J19").Value = 0.3659192825112108
$ws.Range("K19").Value = 0.1201793721973094
$ws.Range("M19").Value = 0.02331838565022422
$ws.Range("N19").Value = 0.0008968609865470852
$ws.Range("O19").Value = 0.06636771300448431
$ws.Range("S19").Value = 0.1058295964125561

Write-Host "Applied 108 value updates to Sheet1"